# Infrared control.pptx - "Add files via upload"
#
# Slide 4, shape "文本框 7" (the paragraph that explains how to download the
# package): the old YB_IR repo URL is replaced by the new
# YahboomTechnology/Yahboom_IR one, and the shape grows a bit taller to fit
# the (slightly longer) wrapped text.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$shp = $s.Shapes.Item(9)              # "文本框 7"
$tr = $shp.TextFrame.TextRange

$oldUrl = "https://github.com/lzty634158/YB_IR "
$newUrl = "https://github.com/YahboomTechnology/Yahboom_IR"

$fullText = $tr.Text
$matchIndex = $fullText.IndexOf($oldUrl)
if ($matchIndex -lt 0) {
    throw "Could not find the expected URL run text in the shape."
}
$startPos = $matchIndex + 1   # TextRange positions are 1-based

# Replace the old URL (which previously included a trailing space as part of
# the same run) with the new URL text, with no trailing space this time.
$urlRange = $tr.Characters($startPos, $oldUrl.Length)
$urlRange.Text = $newUrl

# Re-fetch the range for the freshly written text (the old range object
# keeps referring to the old length) and append a single space after it as
# its own run, so it keeps matching the red URL formatting rather than the
# formatting of the following "to get the package." run.
$newUrlRange = $tr.Characters($startPos, $newUrl.Length)
[void]$newUrlRange.InsertAfter(" ")

$spaceRange = $tr.Characters($startPos + $newUrl.Length, 1)
$spaceRange.Font.Color.RGB = 255   # RGB(255,0,0) -> stored as BGR 0x0000FF -> FF0000

# The text box auto-fits its height to the text; after the edit it grows
# from 175.2pt to 204.0pt (2225040 EMU -> 2590800 EMU).
$shp.Height = 204.0
